$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Docentes responsaveis:" row (old row 12) used to be followed by a
# separate row (old row 13) holding only the teacher's name in B/C with no
# label in A and no custom row height. That extra row is removed entirely
# (its row 13 <row> element goes away) and everything below shifts up by one
# row. Deleting the row reproduces exactly that shift, including updating
# the sheet dimension from A1:C25 to A1:C24.
$ws.Rows.Item(13).Delete()

# After the shift, re-point a handful of cells (which keep their original
# row heights/styles) at their new text content.
$ws.Range("B10:C10").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("B13:C13").Value = "Semestral"
$ws.Range("B15:C15").Value = "01/01/2022"
$ws.Range("B18:C18").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("B19:C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("B20:C20").Value = "Média ponderada de provas  e atividades."
$ws.Range("B21:C21").Value = "1 (uma) prova escrita"
